$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edit -------------------------------------------------------
# Insert a new row above the existing row 2 ("Phone number" data row),
# pushing the current value (9150786843) down to row 3, then fill the
# newly created row 2 with the new phone number.
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = 9790581357

# --- View / selection --------------------------------------------------
# The active selection in the saved sheet view moves to C5.
[void]$ws.Range("C5").Select()

# --- Outline depth bump --------------------------------------------------
# The sheet's stored max row-outline level goes from 1 to 2. Group a row
# outside the used range to two levels deep, then remove that row again so
# no stray grouped row is left behind in the data - this leaves the bumped
# outlineLevelRow recorded on sheetFormatPr, matching the target state.
$ws.Rows.Item(10).Group()
$ws.Rows.Item(10).Group()
$ws.Rows.Item(10).Delete()
